# Applies the weekly price update for Hortaliza - Jengibre (Terminal Hortofrutícola Agro Chillán).
# The data rows (2-20) were re-shuffled: dates in column D and the associated
# Volumen (J), Precio minimo/maximo/promedio (K/L/M) and Precio $/Kg (P) values
# were rearranged across rows. The other columns (A,B,C,E,F,G,H,I,N,O,Q,R) are
# unchanged because they hold the same constant values in every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "11/15/2022"
$ws.Range("J2").Value = 30
$ws.Range("L2").Value = 17000
$ws.Range("M2").Value = 17000
$ws.Range("P2").Value = 1308
# Row 3
$ws.Range("D3").Value = "10/25/2022"
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("P3").Value = 1000
# Row 4
$ws.Range("D4").Value = "08/31/2022"
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12500
$ws.Range("P4").Value = 962
# Row 5
$ws.Range("D5").Value = "12/20/2022"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("P5").Value = 1385
# Row 6
$ws.Range("D6").Value = "11/29/2022"
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 18000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 18000
$ws.Range("P6").Value = 1385
# Row 7
$ws.Range("D7").Value = "10/07/2022"
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 18000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 18000
$ws.Range("P7").Value = 1385
# Row 8
$ws.Range("D8").Value = "10/05/2022"
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 1192
# Row 9
$ws.Range("D9").Value = "08/17/2023"
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 20000
$ws.Range("P9").Value = 1538
# Row 10
$ws.Range("D10").Value = "11/03/2022"
# Row 11
$ws.Range("D11").Value = "09/28/2022"
$ws.Range("J11").Value = 60
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 1346
# Row 12
$ws.Range("D12").Value = "09/06/2022"
$ws.Range("J12").Value = 40
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 12500
$ws.Range("P12").Value = 962
# Row 13
$ws.Range("D13").Value = "08/24/2022"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12500
$ws.Range("P13").Value = 962
# Row 14
$ws.Range("D14").Value = "12/27/2022"
$ws.Range("K14").Value = 17000
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 17000
$ws.Range("P14").Value = 1308
# Row 15
$ws.Range("D15").Value = "01/04/2023"
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 17000
$ws.Range("P15").Value = 1308
# Row 16
$ws.Range("D16").Value = "01/17/2023"
# Row 17
$ws.Range("D17").Value = "11/09/2022"
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 17000
$ws.Range("M17").Value = 17000
$ws.Range("P17").Value = 1308
# Row 18
$ws.Range("D18").Value = "11/30/2022"
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 18000
$ws.Range("P18").Value = 1385
# Row 19
$ws.Range("D19").Value = "10/12/2022"
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 18000
$ws.Range("P19").Value = 1385
# Row 20
$ws.Range("D20").Value = "02/02/2023"
$ws.Range("K20").Value = 19000
$ws.Range("L20").Value = 19000
$ws.Range("M20").Value = 19000
$ws.Range("P20").Value = 1462
